# Database design.xlsx — "Added www interface and cleaning up api"
#
# Adds a new "REST" section header row (bold) plus a row of column labels
# (mirroring the existing ProcessTable-style mini tables) to the bottom of
# the "Log" sheet, and leaves the selection parked on the newly-added range.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Log")

# New bold section title in A25 ("REST"), matching the style used for the
# other section headers on this sheet (bold Calibri 11) but without the
# fill/border those use — just bold text on the default style.
$ws.Range("A25").Value = "REST"
$ws.Range("A25").Font.Bold = $true

# New header row describing the REST log columns, reusing existing labels.
$ws.Range("B25").Value = "PCId"
$ws.Range("C25").Value = "LogCode"
$ws.Range("D25").Value = "Environment"
$ws.Range("E25").Value = "ScheduleId"
$ws.Range("F25").Value = "Message"

# Park the selection on the cell the author last had selected.
[void]$ws.Range("K19").Select()
